$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
foreach ($col in 1..21) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# Add table over the data range
$rng = $ws.Range("A1:U55")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze top row (freeze panes)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
